$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (cohort_year=2021, period_index=4): num_customers 63 -> 64, retention_rate recomputed
$ws.Range("C12").Value = 64
$ws.Range("E12").Value = 0.1534772182254197

# Row 19 (cohort_year=2023, period_index=2): num_customers 60 -> 61, retention_rate recomputed
$ws.Range("C19").Value = 61
$ws.Range("E19").Value = 0.4919354838709677

# Row 22 (cohort_year=2025, period_index=0): num_customers 55 -> 56, cohort_size 55 -> 56
$ws.Range("C22").Value = 56
$ws.Range("D22").Value = 56
